$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.302.91"
$ws.Range("D3").Value = "1.867.72"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.26"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4722"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2867"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.64"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.06"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.863.90"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7154"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.113"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.02"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "30.294.98"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007449"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "2.110.18"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.234"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.23"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.963"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.874"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09596"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.313"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.478"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.198"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.107"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04783"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6831"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01878"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.841"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.38"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.201"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4185"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.582"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.985"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.95"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "891.85"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05739"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.03%  "
